# "Generate Report for Handoff"
#
# The handoff batch for yuwzho's localization job was re-run: a new guid
# (13a2c77f-c290-4c81-9d2e-04bf70b9f435) replaces the old one
# (2fcc8788-5d5b-48bc-85f3-49c4cf484ebc), the xlf checksum segment changes
# from cff0f82b4abf232c58c2d6c276f7d3dd220e610d to
# c7e137a5d06ee2484d96da06ca7842d4c07b259d, and the handoff timestamps move
# forward a little over a minute. Update the three report sheets
# (Overview, zh-cn, de-de) accordingly, keeping each hyperlink pointed at
# its original target URL but with refreshed display text.

$wb = $excel.ActiveWorkbook

$oldGuid = "2fcc8788-5d5b-48bc-85f3-49c4cf484ebc"
$newGuid = "13a2c77f-c290-4c81-9d2e-04bf70b9f435"
$oldHash = "cff0f82b4abf232c58c2d6c276f7d3dd220e610d"
$newHash = "c7e137a5d06ee2484d96da06ca7842d4c07b259d"

$newMdName    = "$newGuid.md"
$newZhXlfName = "$newGuid.$newHash.zh-cn.xlf"
$newDeXlfName = "$newGuid.$newHash.de-de.xlf"

$newOverviewDate = "2016-03-30 10:33:37"
$newZhDate       = "2016-03-30 10:33:27"
$newDeDate       = "2016-03-30 10:33:37"

# NOTE: `Range.Hyperlinks.Delete()` in this host clears EVERY hyperlink on
# the parent worksheet (it is not scoped to the range), and
# `Worksheet.Hyperlinks.Add` always appends a brand-new entry rather than
# replacing one in place, and it also forces a generic blue/underline
# "Hyperlink" font that doesn't match this workbook's existing custom
# HyperLink cell style (underline, Calibri 11, RGB 6495ED). So, per sheet:
# wipe all hyperlinks once, update the cell text, re-add every hyperlink
# that sheet needs (same target address as before, refreshed display
# text), then restore each linked cell's original font so its appearance
# (and resolved style) stays exactly as it was before the edit.

function Set-HyperlinkFont($range) {
    $range.Font.Name = "Calibri"
    $range.Font.Size = 11
    $range.Font.Underline = $true
    $range.Font.Color = 15570276   # RGB(0xED,0x95,0x64) == hex 6495ED
}

# ---- Sheet "Overview" ----
$wsOverview = $wb.Worksheets.Item("Overview")
$overviewMdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/783795447f421ce81910e6ba6dc1b79ff7bee619/e2e/$oldGuid.md"

$wsOverview.Range("A2").Hyperlinks.Delete()
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("D2").Value = $newOverviewDate
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $overviewMdAddress, $null, $null, $newMdName) | Out-Null
Set-HyperlinkFont $wsOverview.Range("A2")

# ---- Sheet "zh-cn" ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$zhMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/783795447f421ce81910e6ba6dc1b79ff7bee619/e2e/$oldGuid.md"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/3fee678beacd499990f1e53d6c9d04a7e0d3a647/ol-handoff/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/ht/$oldGuid.$oldHash.zh-cn.xlf"

$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newMdName
$wsZh.Range("D2").Value = $newZhXlfName
$wsZh.Range("E2").Value = $newZhDate
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $zhMdAddress, $null, $null, $newMdName) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("D2"), $zhXlfAddress, $null, $null, $newZhXlfName) | Out-Null
Set-HyperlinkFont $wsZh.Range("A2")
Set-HyperlinkFont $wsZh.Range("D2")

# ---- Sheet "de-de" ----
$wsDe = $wb.Worksheets.Item("de-de")
$deMdAddress  = "https://github.com/OpenLocalizationTest/oltest/blob/783795447f421ce81910e6ba6dc1b79ff7bee619/e2e/$oldGuid.md"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/74b32dd033eca1937702de351d9d1f4781937738/ol-handoff/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/ht/$oldGuid.$oldHash.de-de.xlf"

$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newMdName
$wsDe.Range("D2").Value = $newDeXlfName
$wsDe.Range("E2").Value = $newDeDate
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $deMdAddress, $null, $null, $newMdName) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("D2"), $deXlfAddress, $null, $null, $newDeXlfName) | Out-Null
Set-HyperlinkFont $wsDe.Range("A2")
Set-HyperlinkFont $wsDe.Range("D2")

Write-Output "Handoff report regenerated."
